$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 21:24"

# 2. Fix ordering of country names that moved in the shared string table
#    (Cabo Verde moves to before Cuba; Botsuana moves to before Malta)
$ws.Range("A121").Value = "Cabo Verde"
$ws.Range("A122").Value = "Cuba"
$ws.Range("A123").Value = "Surinam"

$ws.Range("A146").Value = "Botsuana"
$ws.Range("A147").Value = "Malta"
$ws.Range("A148").Value = "Sierra Leona"

# 3. Update the updated statistics for all affected rows (Casos totales,
#    Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6475505
$ws.Range("C4").Value = 15255
$ws.Range("D4").Value = 3735798
$ws.Range("E4").Value = 2546308
$ws.Range("G4").Value = 151
$ws.Range("H4").Value = 193399

# Row 5: India
$ws.Range("B5").Value = 4276777
$ws.Range("C5").Value = 74215
$ws.Range("D5").Value = 3320947
$ws.Range("E5").Value = 883021
$ws.Range("G5").Value = 1122
$ws.Range("H5").Value = 72809

# Row 17: Francia
$ws.Range("D17").Value = 87836
$ws.Range("E17").Value = 210418
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 30726

# Row 32: Ecuador
$ws.Range("B32").Value = 110092
$ws.Range("C32").Value = 308
$ws.Range("E32").Value = 8274
$ws.Range("G32").Value = 52
$ws.Range("H32").Value = 10576

# Row 91: Grecia
$ws.Range("B91").Value = 11663
$ws.Range("C91").Value = 139
$ws.Range("E91").Value = 7570

# Row 92: Noruega
$ws.Range("B92").Value = 11462
$ws.Range("C92").Value = 74
$ws.Range("E92").Value = 1850

# Row 121: Cabo Verde (new data)
$ws.Range("B121").Value = 4358
$ws.Range("C121").Value = 28
$ws.Range("D121").Value = 3790
$ws.Range("E121").Value = 526
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 42

# Row 122: Cuba
$ws.Range("B122").Value = 4352
$ws.Range("C122").Value = 43
$ws.Range("D122").Value = 3642
$ws.Range("E122").Value = 608
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 102

# Row 123: Surinam
$ws.Range("B123").Value = 4346
$ws.Range("D123").Value = 3494
$ws.Range("E123").Value = 767
$ws.Range("H123").Value = 85

# Row 134: Angola
$ws.Range("B134").Value = 2981
$ws.Range("C134").Value = 16
$ws.Range("D134").Value = 1215
$ws.Range("E134").Value = 1646
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 120

# Row 146: Botsuana (new data)
$ws.Range("B146").Value = 2126
$ws.Range("C146").Value = 124
$ws.Range("D146").Value = 493
$ws.Range("E146").Value = 1624
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 9

# Row 147: Malta
$ws.Range("B147").Value = 2076
$ws.Range("C147").Value = 37
$ws.Range("D147").Value = 1690
$ws.Range("H147").Value = 14

# Row 148: Sierra Leona
$ws.Range("B148").Value = 2054
$ws.Range("D148").Value = 1611
$ws.Range("E148").Value = 372
$ws.Range("H148").Value = 71

# Row 149: Yemen
$ws.Range("B149").Value = 1989
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 1201
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 573

# Row 191: Monaco
$ws.Range("B191").Value = 153
$ws.Range("C191").Value = 6
$ws.Range("E191").Value = 50

# Row 193: Seychelles
$ws.Range("B193").Value = 137
$ws.Range("C193").Value = 1
$ws.Range("E193").Value = 10
